$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and name it "test_data3"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "test_data3"

# Populate the table data (mirrors the layout used in the other test_data sheets)
$newSheet.Range("C14").Value = "Sweden_Pay_Now_Direct_debit"

$newSheet.Range("C15").Value = "column_name1"
$newSheet.Range("D15").Value = "value1"

$newSheet.Range("C16").Value = "column_name2"
$newSheet.Range("D16").Value = "value2"

$newSheet.Range("C17").Value = "column_name3"
$newSheet.Range("D17").Value = "value3"

$newSheet.Range("C18").Value = "column_name4"
$newSheet.Range("D18").Value = "value4"
